$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("Q5").Value = 2.3
$ws.Range("R5").Value = 1.6
$ws.Range("N7").Value = 13.8
$ws.Range("P7").Value = 4.2
$ws.Range("S7").Value = 1.29
$ws.Range("T7").Value = 3.42
$ws.Range("G8").Value = 3.8
$ws.Range("J8").Value = 4.2
$ws.Range("K8").Value = 2.12
$ws.Range("U8").Value = 1.7
$ws.Range("V8").Value = 1.93
$ws.Range("W8").Value = 11.5
$ws.Range("X8").Value = 21
$ws.Range("Y8").Value = 12.5
$ws.Range("AD8").Value = 6.7
$ws.Range("AG8").Value = 7.6
$ws.Range("AH8").Value = 9.25
$ws.Range("AJ8").Value = 16
$ws.Range("AK8").Value = 14.5
$ws.Range("AN8").Value = 5.6
$ws.Range("AO8").Value = 21
$ws.Range("AQ8").Value = 110
$ws.Range("AT8").Value = 2.75
$ws.Range("AY8").Value = 17.5
$ws.Range("M15").Value = 1.08
$ws.Range("N15").Value = 8
$ws.Range("O15").Value = 1.4
$ws.Range("P15").Value = 2.75
$ws.Range("Q15").Value = 2.35
$ws.Range("R15").Value = 1.57
$ws.Range("W15").Value = 7
$ws.Range("AD15").Value = 5.5
$ws.Range("AF15").Value = 51
$ws.Range("AM15").Value = 800
$ws.Range("AU15").Value = 8.5
$ws.Range("AW15").Value = 4.75
$ws.Range("M17").Value = 1.04
$ws.Range("N17").Value = 13
$ws.Range("Q20").Value = 2.25
$ws.Range("R20").Value = 1.62
$ws.Range("Q21").Value = 2.4
$ws.Range("R21").Value = 1.53
$ws.Range("J26").Value = 2.7
$ws.Range("L26").Value = 3.5
$ws.Range("S26").Value = 1.36
$ws.Range("T26").Value = 2.92
$ws.Range("X26").Value = 11
$ws.Range("AB26").Value = 24
$ws.Range("AG26").Value = 11.25
$ws.Range("AH26").Value = 17.5
$ws.Range("AJ26").Value = 37
$ws.Range("AO26").Value = 11
$ws.Range("AP26").Value = 18
$ws.Range("AR26").Value = 70
$ws.Range("AT26").Value = 2.92
$ws.Range("AX26").Value = 16
